# Add three new Zika outbreak records (rows 46-48) to Sheet1:
#   - Colombia / Norte de Santander (row 46)
#   - El Salvador / Usultan (row 47)
#   - El Salvador / San Salvador - Tonacatepeque (row 48)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 46 - Colombia, Norte de Santander
$ws.Range("A46").Value = "Colombia"
$ws.Range("D46").Value = "-"
$ws.Range("E46").Value = "1.33 million (2013)"
$ws.Range("F46").Value = "-"
$ws.Range("G46").Value = "Norte de Santander"
$ws.Range("H46").Value = "2015-18-11"
$ws.Range("I46").Value = "Zika"
$ws.Range("J46").Value = "44 cases confimed in Colombia (Cucuta, El Zulia, Puerto Santander, San Cayetano, Tibu, & Villa Del Rosario) "
$ws.Range("K46").Value = "http://diariodelcauca.com.co/noticias/nacional/ascienden-393-los-casos-de-contagio-por-virus-zika-en-colo-170494"

# Row 47 - El Salvador, Usultan (Southeast of region)
$ws.Range("A47").Value = "El Salvador"
$ws.Range("E47").Value = "366 000"
$ws.Range("F47").Value = "Southeast of region"
$ws.Range("G47").Value = "Usultan"
$ws.Range("H47").Value = "2015-25-11"
$ws.Range("I47").Value = "Zika"
$ws.Range("J47").Value = "33 and 28-year old woman (first three cases in El Salvador)"
$ws.Range("K47").Value = "http://www.telemetro.com/actualidad/salud/Salvador-confirma-presencia-virus-pais_0_865414161.html"

# Row 48 - El Salvador, San Salvador - Tonacatepeque
$ws.Range("A48").Value = "El Salvador"
$ws.Range("B48").Value = 13.791
$ws.Range("C48").Value = -89.286000000000001
$ws.Range("D48").Value = 603
$ws.Range("E48").Value = "90 896 (2007)"
$ws.Range("F48").Value = "Tonacatepeque"
$ws.Range("G48").Value = "San Salvador"
$ws.Range("H48").Value = "2015-25-11"
$ws.Range("I48").Value = "Zika"
$ws.Range("J48").Value = "33-year old (first three cases in El Salvador)"
$ws.Range("K48").Value = "http://www.telemetro.com/actualidad/salud/Salvador-confirma-presencia-virus-pais_0_865414161.html"

# Match the author's final view state: scrolled down a bit, with C49 selected.
$ws.Activate()
$ws.Range("C49").Select()
